$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (before) values of rows 14 and 15
$row14 = @{
    A = $ws.Cells.Item(14,1).Value()
    B = $ws.Cells.Item(14,2).Value()
    C = $ws.Cells.Item(14,3).Value()
    E = $ws.Cells.Item(14,5).Value()
    F = $ws.Cells.Item(14,6).Value()
    G = $ws.Cells.Item(14,7).Value()
}

$row15 = @{
    A = $ws.Cells.Item(15,1).Value()
    B = $ws.Cells.Item(15,2).Value()
    C = $ws.Cells.Item(15,3).Value()
    E = $ws.Cells.Item(15,5).Value()
    F = $ws.Cells.Item(15,6).Value()
    G = $ws.Cells.Item(15,7).Value()
}

# Row 14 becomes the old row 15 content (Alloy modelization task), keep shared formula in D14
$ws.Cells.Item(14,1).Value = $row15.A
$ws.Cells.Item(14,2).Value = $row15.B
$ws.Cells.Item(14,3).Value = $row15.C
$ws.Cells.Item(14,5).Value = $row15.E
$ws.Cells.Item(14,6).Value = $row15.F
$ws.Cells.Item(14,7).Value = $row15.G

# Row 15 becomes the old row 14 content (Document revision task)
$ws.Cells.Item(15,1).Value = $row14.A
$ws.Cells.Item(15,2).Value = $row14.B
$ws.Cells.Item(15,3).Value = $row14.C
$ws.Cells.Item(15,5).Value = $row14.E
$ws.Cells.Item(15,6).Value = $row14.F
$ws.Cells.Item(15,7).Value = $row14.G

# D15 formula becomes a standalone (non-shared) formula instead of part of the shared group
$ws.Cells.Item(15,4).Formula = '=CONCATENATE(NETWORKDAYS(E15,F15),"g")'

# Update selection to reflect G15 as active cell
$ws.Range("G15").Select() | Out-Null
